# Applies a cyclic rotation of the species-observation rows 4, 5 and 6:
#   new row 4 <- old row 6
#   new row 5 <- old row 4
#   new row 6 <- old row 5
# Only columns A,B,D,E,F,G,H,I,J,Q,R actually differ between the three
# rows (every other column already holds identical values across the
# trio), so only those are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    $data = @{}
    $data.A = $ws.Range("A$row").Value2
    $data.B = $ws.Range("B$row").Value2
    $data.D = $ws.Range("D$row").Value2
    $data.E = $ws.Range("E$row").Value2
    $data.F = $ws.Range("F$row").Value2
    $data.G = $ws.Range("G$row").Value2
    $data.H = $ws.Range("H$row").Value2
    $data.I = $ws.Range("I$row").Text
    $data.J = $ws.Range("J$row").Text
    $data.Q = $ws.Range("Q$row").Value2
    $data.R = $ws.Range("R$row").Value2
    return $data
}

function Set-RowData($row, $data) {
    $ws.Range("A$row").Value = $data.A
    $ws.Range("B$row").Value = $data.B
    $ws.Range("D$row").Value = $data.D
    $ws.Range("E$row").Value = $data.E
    $ws.Range("F$row").Value = $data.F
    $ws.Range("G$row").Value = $data.G
    $ws.Range("H$row").Value = $data.H

    # Column I holds numeric-looking text ("5", "1", ...); only touch it
    # when the value actually needs to change, so rows that keep the same
    # digit are left byte-for-byte alone. Briefly force a text number
    # format so the numeric-looking string is not re-interpreted as a
    # number, then restore the default "Normal" style so no stray
    # formatting is left behind.
    if ($ws.Range("I$row").Text -ne $data.I) {
        $ws.Range("I$row").NumberFormat = "@"
        $ws.Range("I$row").Value = $data.I
        $ws.Range("I$row").Style = "Normal"
    }

    # Column J ("m²", "plantor/tuvor", ...) is never numeric-looking, so it
    # round-trips as text without any extra formatting; again, only write
    # it when it actually changes.
    if ($ws.Range("J$row").Text -ne $data.J) {
        $ws.Range("J$row").Value = $data.J
    }

    $ws.Range("Q$row").Value = $data.Q
    $ws.Range("R$row").Value = $data.R
}

$row4 = Get-RowData 4
$row5 = Get-RowData 5
$row6 = Get-RowData 6

Set-RowData 4 $row6
Set-RowData 5 $row4
Set-RowData 6 $row5
